# Update well labels from E1..E12 to A1..A12 in column A (rows 2-13)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($i = 1; $i -le 12; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = "A$i"
}

# Update the selection to match the new edit state
$ws.Range("A2:A13").Select()
